$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price cells whose new values would
# otherwise be auto-detected as numbers by Excel (they must stay text,
# matching the original inline-string cells).
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D15","D16","D17","D18","D19","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.055.49"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.884.59"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "307.32"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "0.5170"
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("D8").Value = "0.3724"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("D9").Value = "0.07207"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "0.9035"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "20.94"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07602"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "95.17"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.855.74"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "5.267"
$ws.Range("E15").Value = "  +0.62%  "
$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "0.000008502"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "14.29"
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "0.9981"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "27.092.67"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "5.042"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "2.124.93"
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").Value = "10.47"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").Value = "6.460"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").Value = "145.82"
$ws.Range("D26").Value = "1.789"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").Value = "18.01"
$ws.Range("E27").Value = "  +1.01%  "
$ws.Range("D28").Value = "2.121"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").Value = "114.64"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "4.927"
$ws.Range("E30").Value = "  +5.48%  "
$ws.Range("D31").Value = "4.784"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").Value = "0.09212"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "0.05043"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").Value = "0.7629"
$ws.Range("E34").Value = "  +4.59%  "
$ws.Range("D35").Value = "1.188"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").Value = "3.015"
$ws.Range("E36").Value = "  -1.74%  "
$ws.Range("D37").Value = "3.272"
$ws.Range("E37").Value = "  +2.73%  "
$ws.Range("D38").Value = "2.524"
$ws.Range("E38").Value = "  +3.24%  "
$ws.Range("D39").Value = "0.5606"
$ws.Range("E39").Value = "  +6.33%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "1.076"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "6.595"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "8.900"
$ws.Range("E43").Value = "  +5.70%  "
$ws.Range("D44").Value = "118.05"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").Value = "0.1506"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").Value = "0.4793"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "0.9983"
$ws.Range("E47").Value = "  -0.21%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "10.14"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "1.574"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "37.13"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "63.58"
$ws.Range("E51").Value = "  +1.17%  "
